# Atualizei dados da ADD
# Insert a new daily record (Dia=11, total_venda=19798.81, Mes=6/Junho, Ano=2025, Periodo=06/2025)
# at the top of the June block (sheet row 9), shifting the existing rows down by one.
# Also correct the total_venda value for 02/04/2025 (new row 32) from 48732.41 to 48690.41.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 9 (start of the "June" data block)
$ws.Rows.Item(9).Insert()

# Fill in the new row with the added daily record
$ws.Range("A9").Value2 = 11
$ws.Range("B9").Value2 = 19798.81
$ws.Range("C9").Value2 = 6
$ws.Range("D9").Value2 = 2025
$ws.Range("E9").Value2 = "06/2025"

# Correct the value that now sits at row 32 (02/04/2025)
$ws.Range("B32").Value2 = 48690.41

Write-Output "edit applied"
